$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: swap paired rows (columns B:AC), these represent matches whose
# home/away ordering in the source feed changed between extracts. Column A (row id) is untouched.
$rangeA = $ws.Range("B2:AC2")
$rangeB = $ws.Range("B3:AC3")
$valsA = $rangeA.Value()
$valsB = $rangeB.Value()
$rangeA.Value = $valsB
$rangeB.Value = $valsA

$rangeA = $ws.Range("B19:AC19")
$rangeB = $ws.Range("B20:AC20")
$valsA = $rangeA.Value()
$valsB = $rangeB.Value()
$rangeA.Value = $valsB
$rangeB.Value = $valsA

$rangeA = $ws.Range("B31:AC31")
$rangeB = $ws.Range("B32:AC32")
$valsA = $rangeA.Value()
$valsB = $rangeB.Value()
$rangeA.Value = $valsB
$rangeB.Value = $valsA

$rangeA = $ws.Range("B50:AC50")
$rangeB = $ws.Range("B51:AC51")
$valsA = $rangeA.Value()
$valsB = $rangeB.Value()
$rangeA.Value = $valsB
$rangeB.Value = $valsA

$rangeA = $ws.Range("B90:AC90")
$rangeB = $ws.Range("B91:AC91")
$valsA = $rangeA.Value()
$valsB = $rangeB.Value()
$rangeA.Value = $valsB
$rangeB.Value = $valsA

$rangeA = $ws.Range("B124:AC124")
$rangeB = $ws.Range("B125:AC125")
$valsA = $rangeA.Value()
$valsB = $rangeB.Value()
$rangeA.Value = $valsB
$rangeB.Value = $valsA

$rangeA = $ws.Range("B140:AC140")
$rangeB = $ws.Range("B141:AC141")
$valsA = $rangeA.Value()
$valsB = $rangeB.Value()
$rangeA.Value = $valsB
$rangeB.Value = $valsA

$rangeA = $ws.Range("B167:AC167")
$rangeB = $ws.Range("B168:AC168")
$valsA = $rangeA.Value()
$valsB = $rangeB.Value()
$rangeA.Value = $valsB
$rangeB.Value = $valsA

# --- Step 2: append 8 new upcoming-fixture rows (178-185) at the bottom of the table
# row 178
$ws.Cells.Item(178, 1).Value = 176
$ws.Cells.Item(178, 1).Font.Bold = $true
$ws.Cells.Item(178, 1).HorizontalAlignment = -4108
$ws.Cells.Item(178, 1).VerticalAlignment = -4160
$ws.Cells.Item(178, 1).Borders.LineStyle = 1
$ws.Cells.Item(178, 2).Value = 6979559
$ws.Cells.Item(178, 3).Value = 'Serbia Super Liga'
$ws.Cells.Item(178, 4).Value = 'Serbia Super Liga'
$ws.Cells.Item(178, 5).Value = 45352.416666666664
$ws.Cells.Item(178, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(178, 6).Value = 'FK Radnicki 1923'
$ws.Cells.Item(178, 7).Value = 'IMT Novi Belgrade'
$ws.Cells.Item(178, 11).Value = 1.8
$ws.Cells.Item(178, 12).Value = 3.4
$ws.Cells.Item(178, 13).Value = 4.2
$ws.Cells.Item(178, 14).Value = 1.85
$ws.Cells.Item(178, 15).Value = 3.3
$ws.Cells.Item(178, 16).Value = 4
$ws.Cells.Item(178, 17).Value = -0.5
$ws.Cells.Item(178, 18).Value = 1.9
$ws.Cells.Item(178, 19).Value = 1.9
$ws.Cells.Item(178, 20).Value = 2.5
$ws.Cells.Item(178, 21).Value = 2
$ws.Cells.Item(178, 22).Value = 1.8
$ws.Cells.Item(178, 23).Value = 0
$ws.Cells.Item(178, 24).Value = 0
$ws.Cells.Item(178, 25).Value = 0
$ws.Cells.Item(178, 26).Value = 0
$ws.Cells.Item(178, 27).Value = 0

# row 179
$ws.Cells.Item(179, 1).Value = 177
$ws.Cells.Item(179, 1).Font.Bold = $true
$ws.Cells.Item(179, 1).HorizontalAlignment = -4108
$ws.Cells.Item(179, 1).VerticalAlignment = -4160
$ws.Cells.Item(179, 1).Borders.LineStyle = 1
$ws.Cells.Item(179, 2).Value = 6979562
$ws.Cells.Item(179, 3).Value = 'Serbia Super Liga'
$ws.Cells.Item(179, 4).Value = 'Serbia Super Liga'
$ws.Cells.Item(179, 5).Value = 45352.604166666664
$ws.Cells.Item(179, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(179, 6).Value = 'FK Cukaricki'
$ws.Cells.Item(179, 7).Value = 'FK Napredak'
$ws.Cells.Item(179, 11).Value = 1.5
$ws.Cells.Item(179, 12).Value = 4
$ws.Cells.Item(179, 13).Value = 5.75
$ws.Cells.Item(179, 14).Value = 1.533
$ws.Cells.Item(179, 15).Value = 4
$ws.Cells.Item(179, 16).Value = 5.5
$ws.Cells.Item(179, 17).Value = -1
$ws.Cells.Item(179, 18).Value = 1.925
$ws.Cells.Item(179, 19).Value = 1.875
$ws.Cells.Item(179, 20).Value = 2.5
$ws.Cells.Item(179, 21).Value = 1.875
$ws.Cells.Item(179, 22).Value = 1.925
$ws.Cells.Item(179, 23).Value = 0
$ws.Cells.Item(179, 24).Value = 0
$ws.Cells.Item(179, 25).Value = 0
$ws.Cells.Item(179, 26).Value = 0
$ws.Cells.Item(179, 27).Value = 0

# row 180
$ws.Cells.Item(180, 1).Value = 178
$ws.Cells.Item(180, 1).Font.Bold = $true
$ws.Cells.Item(180, 1).HorizontalAlignment = -4108
$ws.Cells.Item(180, 1).VerticalAlignment = -4160
$ws.Cells.Item(180, 1).Borders.LineStyle = 1
$ws.Cells.Item(180, 2).Value = 6979558
$ws.Cells.Item(180, 3).Value = 'Serbia Super Liga'
$ws.Cells.Item(180, 4).Value = 'Serbia Super Liga'
$ws.Cells.Item(180, 5).Value = 45353.416666666664
$ws.Cells.Item(180, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(180, 6).Value = 'Javor Ivanjica'
$ws.Cells.Item(180, 7).Value = 'Vojvodina'
$ws.Cells.Item(180, 11).Value = 3
$ws.Cells.Item(180, 12).Value = 3.3
$ws.Cells.Item(180, 13).Value = 2.2
$ws.Cells.Item(180, 14).Value = 3.1
$ws.Cells.Item(180, 15).Value = 3.3
$ws.Cells.Item(180, 16).Value = 2.1
$ws.Cells.Item(180, 17).Value = 0.25
$ws.Cells.Item(180, 18).Value = 1.95
$ws.Cells.Item(180, 19).Value = 1.85
$ws.Cells.Item(180, 20).Value = 2.5
$ws.Cells.Item(180, 21).Value = 1.975
$ws.Cells.Item(180, 22).Value = 1.825
$ws.Cells.Item(180, 23).Value = 0
$ws.Cells.Item(180, 24).Value = 0
$ws.Cells.Item(180, 25).Value = 0
$ws.Cells.Item(180, 26).Value = 0
$ws.Cells.Item(180, 27).Value = 0

# row 181
$ws.Cells.Item(181, 1).Value = 179
$ws.Cells.Item(181, 1).Font.Bold = $true
$ws.Cells.Item(181, 1).HorizontalAlignment = -4108
$ws.Cells.Item(181, 1).VerticalAlignment = -4160
$ws.Cells.Item(181, 1).Borders.LineStyle = 1
$ws.Cells.Item(181, 2).Value = 6978757
$ws.Cells.Item(181, 3).Value = 'Serbia Super Liga'
$ws.Cells.Item(181, 4).Value = 'Serbia Super Liga'
$ws.Cells.Item(181, 5).Value = 45353.5
$ws.Cells.Item(181, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(181, 6).Value = 'FK Backa Topola'
$ws.Cells.Item(181, 7).Value = 'Crvena Zvezda'
$ws.Cells.Item(181, 11).Value = 6
$ws.Cells.Item(181, 12).Value = 4.75
$ws.Cells.Item(181, 13).Value = 1.4
$ws.Cells.Item(181, 14).Value = 5.5
$ws.Cells.Item(181, 15).Value = 4.75
$ws.Cells.Item(181, 16).Value = 1.444
$ws.Cells.Item(181, 17).Value = 1.25
$ws.Cells.Item(181, 18).Value = 1.875
$ws.Cells.Item(181, 19).Value = 1.925
$ws.Cells.Item(181, 20).Value = 3
$ws.Cells.Item(181, 21).Value = 1.85
$ws.Cells.Item(181, 22).Value = 1.95
$ws.Cells.Item(181, 23).Value = 0
$ws.Cells.Item(181, 24).Value = 0
$ws.Cells.Item(181, 25).Value = 0
$ws.Cells.Item(181, 26).Value = 0
$ws.Cells.Item(181, 27).Value = 0

# row 182
$ws.Cells.Item(182, 1).Value = 180
$ws.Cells.Item(182, 1).Font.Bold = $true
$ws.Cells.Item(182, 1).HorizontalAlignment = -4108
$ws.Cells.Item(182, 1).VerticalAlignment = -4160
$ws.Cells.Item(182, 1).Borders.LineStyle = 1
$ws.Cells.Item(182, 2).Value = 6979561
$ws.Cells.Item(182, 3).Value = 'Serbia Super Liga'
$ws.Cells.Item(182, 4).Value = 'Serbia Super Liga'
$ws.Cells.Item(182, 5).Value = 45353.604166666664
$ws.Cells.Item(182, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(182, 6).Value = 'Partizan Belgrade'
$ws.Cells.Item(182, 7).Value = 'FK Zeleznicar Pancevo'
$ws.Cells.Item(182, 11).Value = 1.181
$ws.Cells.Item(182, 12).Value = 6.5
$ws.Cells.Item(182, 13).Value = 11
$ws.Cells.Item(182, 14).Value = 1.181
$ws.Cells.Item(182, 15).Value = 6.5
$ws.Cells.Item(182, 16).Value = 11
$ws.Cells.Item(182, 17).Value = -2
$ws.Cells.Item(182, 18).Value = 1.875
$ws.Cells.Item(182, 19).Value = 1.925
$ws.Cells.Item(182, 20).Value = 3.25
$ws.Cells.Item(182, 21).Value = 1.825
$ws.Cells.Item(182, 22).Value = 1.975
$ws.Cells.Item(182, 23).Value = 0
$ws.Cells.Item(182, 24).Value = 0
$ws.Cells.Item(182, 25).Value = 0
$ws.Cells.Item(182, 26).Value = 0
$ws.Cells.Item(182, 27).Value = 0

# row 183
$ws.Cells.Item(183, 1).Value = 181
$ws.Cells.Item(183, 1).Font.Bold = $true
$ws.Cells.Item(183, 1).HorizontalAlignment = -4108
$ws.Cells.Item(183, 1).VerticalAlignment = -4160
$ws.Cells.Item(183, 1).Borders.LineStyle = 1
$ws.Cells.Item(183, 2).Value = 6979563
$ws.Cells.Item(183, 3).Value = 'Serbia Super Liga'
$ws.Cells.Item(183, 4).Value = 'Serbia Super Liga'
$ws.Cells.Item(183, 5).Value = 45354.4375
$ws.Cells.Item(183, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(183, 6).Value = 'FK Novi Pazar'
$ws.Cells.Item(183, 7).Value = 'Mladost Lucani'
$ws.Cells.Item(183, 11).Value = 1.909
$ws.Cells.Item(183, 12).Value = 3.3
$ws.Cells.Item(183, 13).Value = 3.75
$ws.Cells.Item(183, 14).Value = 1.85
$ws.Cells.Item(183, 15).Value = 3.3
$ws.Cells.Item(183, 16).Value = 3.8
$ws.Cells.Item(183, 17).Value = -0.5
$ws.Cells.Item(183, 18).Value = 1.9
$ws.Cells.Item(183, 19).Value = 1.9
$ws.Cells.Item(183, 20).Value = 2.5
$ws.Cells.Item(183, 21).Value = 1.9
$ws.Cells.Item(183, 22).Value = 1.9
$ws.Cells.Item(183, 23).Value = 0
$ws.Cells.Item(183, 24).Value = 0
$ws.Cells.Item(183, 25).Value = 0
$ws.Cells.Item(183, 26).Value = 0
$ws.Cells.Item(183, 27).Value = 0

# row 184
$ws.Cells.Item(184, 1).Value = 182
$ws.Cells.Item(184, 1).Font.Bold = $true
$ws.Cells.Item(184, 1).HorizontalAlignment = -4108
$ws.Cells.Item(184, 1).VerticalAlignment = -4160
$ws.Cells.Item(184, 1).Borders.LineStyle = 1
$ws.Cells.Item(184, 2).Value = 6979564
$ws.Cells.Item(184, 3).Value = 'Serbia Super Liga'
$ws.Cells.Item(184, 4).Value = 'Serbia Super Liga'
$ws.Cells.Item(184, 5).Value = 45354.520833333336
$ws.Cells.Item(184, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(184, 6).Value = 'FK Vozdovac'
$ws.Cells.Item(184, 7).Value = 'Spartak Subotica'
$ws.Cells.Item(184, 11).Value = 2.1
$ws.Cells.Item(184, 12).Value = 3.25
$ws.Cells.Item(184, 13).Value = 3.25
$ws.Cells.Item(184, 14).Value = 2.1
$ws.Cells.Item(184, 15).Value = 3.25
$ws.Cells.Item(184, 16).Value = 3.25
$ws.Cells.Item(184, 17).Value = -0.25
$ws.Cells.Item(184, 18).Value = 1.85
$ws.Cells.Item(184, 19).Value = 1.95
$ws.Cells.Item(184, 20).Value = 2.5
$ws.Cells.Item(184, 21).Value = 1.975
$ws.Cells.Item(184, 22).Value = 1.825
$ws.Cells.Item(184, 23).Value = 0
$ws.Cells.Item(184, 24).Value = 0
$ws.Cells.Item(184, 25).Value = 0
$ws.Cells.Item(184, 26).Value = 0
$ws.Cells.Item(184, 27).Value = 0

# row 185
$ws.Cells.Item(185, 1).Value = 183
$ws.Cells.Item(185, 1).Font.Bold = $true
$ws.Cells.Item(185, 1).HorizontalAlignment = -4108
$ws.Cells.Item(185, 1).VerticalAlignment = -4160
$ws.Cells.Item(185, 1).Borders.LineStyle = 1
$ws.Cells.Item(185, 2).Value = 6979560
$ws.Cells.Item(185, 3).Value = 'Serbia Super Liga'
$ws.Cells.Item(185, 4).Value = 'Serbia Super Liga'
$ws.Cells.Item(185, 5).Value = 45354.604166666664
$ws.Cells.Item(185, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(185, 6).Value = 'Radnicki Nis'
$ws.Cells.Item(185, 7).Value = 'FK Radnik Surdulica'
$ws.Cells.Item(185, 11).Value = 1.65
$ws.Cells.Item(185, 12).Value = 3.5
$ws.Cells.Item(185, 13).Value = 5
$ws.Cells.Item(185, 14).Value = 1.95
$ws.Cells.Item(185, 15).Value = 3.2
$ws.Cells.Item(185, 16).Value = 3.75
$ws.Cells.Item(185, 17).Value = -0.5
$ws.Cells.Item(185, 18).Value = 2
$ws.Cells.Item(185, 19).Value = 1.8
$ws.Cells.Item(185, 20).Value = 2.5
$ws.Cells.Item(185, 21).Value = 2
$ws.Cells.Item(185, 22).Value = 1.8
$ws.Cells.Item(185, 23).Value = 0
$ws.Cells.Item(185, 24).Value = 0
$ws.Cells.Item(185, 25).Value = 0
$ws.Cells.Item(185, 26).Value = 0
$ws.Cells.Item(185, 27).Value = 0

